# clean-up of input tables
# Rename the worksheet from "updated" to "Tabelle1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Tabelle1"
